# Update "想去人数" (want-to-go count) figures to the latest generated output.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1847
$wsExpo.Range("F4").Value = 811
$wsExpo.Range("F5").Value = 594
$wsExpo.Range("F6").Value = 226

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1847
$wsAll.Range("F5").Value = 811
$wsAll.Range("F6").Value = 594
$wsAll.Range("F7").Value = 226
